# Scenario 1.xlsx edit script
# - Decrease probability of fetal death before 4 weeks from conception
#   (Phase1!B2:B5 and Phase2!C2:C5), letting the dependent formulas
#   (Phase1!D2:D5, Phase2!E2:E5) recalculate automatically.
# - Clear the stray fill-only format that had been applied to
#   Phase2!D39:D41 and Phase3!G1 (reverts them to the default style).
# - Update the saved view/selection state to match what was left
#   active when the workbook was last saved.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Phase1: rows 2-5, column B (probability of fetal death)
# ---------------------------------------------------------------
$phase1 = $wb.Worksheets.Item("Phase1")
$phase1.Range("B2").Value = 0.1
$phase1.Range("B3").Value = 0.1
$phase1.Range("B4").Value = 0.05
$phase1.Range("B5").Value = 0.05

# ---------------------------------------------------------------
# Phase2: rows 2-5, column C (probability of fetal death)
# ---------------------------------------------------------------
$phase2 = $wb.Worksheets.Item("Phase2")
$phase2.Range("C2").Value = 0.1
$phase2.Range("C3").Value = 0.1
$phase2.Range("C4").Value = 0.05
$phase2.Range("C5").Value = 0.05

# Remove the (unused / stray) fill-applied format from D39:D41 so the
# cells fall back to the default style.
$phase2.Range("D39:D41").ClearFormats()

# ---------------------------------------------------------------
# Phase3: remove the same stray format from G1
# ---------------------------------------------------------------
$phase3 = $wb.Worksheets.Item("Phase3")
$phase3.Range("G1").ClearFormats()

# ---------------------------------------------------------------
# View / selection bookkeeping to match the saved workbook state
# ---------------------------------------------------------------
$phase1.Activate()
$phase1.Range("B2:B5").Select()

$phase2.Range("C2:C5").Select()

$phase5 = $wb.Worksheets.Item("Phase5")
$phase5.Range("K14").Select()

$excel.ActiveWindow.WindowState = -4143
